$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 56 used to be "Lanzarote" (69,137,1564,3) and row 57 used to be
# "La Palma" (69,137,1564,4). Swap the province labels and the Muertes
# (deaths) figures so that La Palma now appears in row 56 and Lanzarote in
# row 57.
$ws.Range("A56").Value = "La Palma"
$ws.Range("E56").Value = 4

$ws.Range("A57").Value = "Lanzarote"
$ws.Range("E57").Value = 3

# Update the "last updated" timestamp string in the title cell (A1).
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 16:52"
